$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '68.299.05'
Set-TextValue "E2" '  +0.62%  '
Set-TextValue "D3" '3.351.97'
Set-TextValue "E3" '  +0.59%  '
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  +0.13%  '
Set-TextValue "D5" '583.77'
Set-TextValue "E5" '  +0.14%  '
Set-TextValue "D6" '177.31'
Set-TextValue "E6" '  +0.54%  '
Set-TextValue "E7" '  -0.05%  '
Set-TextValue "E8" '  +0.12%  '
Set-TextValue "E9" '  +2.94%  '
Set-TextValue "E10" '  +0.82%  '
Set-TextValue "D11" '48.18'
Set-TextValue "E11" '  +5.48%  '
Set-TextValue "E12" '  +1.20%  '
Set-TextValue "D13" '685.29'
Set-TextValue "E13" '  +3.92%  '
Set-TextValue "D14" '3.894.94'
Set-TextValue "E14" '  +0.65%  '
Set-TextValue "D15" '8.41'
Set-TextValue "E15" '  +0.05%  '
Set-TextValue "D16" '68.338.59'
Set-TextValue "E16" '  +0.58%  '
Set-TextValue "E17" '  +1.18%  '
Set-TextValue "D18" '3.354.32'
Set-TextValue "E18" '  +0.67%  '
Set-TextValue "D19" '17.46'
Set-TextValue "E19" '  +0.04%  '
Set-TextValue "D20" '11.19'
Set-TextValue "E20" '  +2.04%  '
Set-TextValue "E21" '  +0.53%  '
Set-TextValue "D22" '5.45'
Set-TextValue "E22" '  +0.30%  '
Set-TextValue "D23" '16.94'
Set-TextValue "E23" '  -0.80%  '
Set-TextValue "D24" '99.87'
Set-TextValue "E24" '  +0.28%  '
Set-TextValue "E25" '  +1.47%  '
Set-TextValue "E27" '  +2.73%  '
Set-TextValue "D28" '32.94'
Set-TextValue "E28" '  -1.88%  '
Set-TextValue "E29" '  +0.57%  '
Set-TextValue "E30" '  -6.61%  '
Set-TextValue "D31" '562.30'
Set-TextValue "E31" '  -5.18%  '
Set-TextValue "D32" '11.07'
Set-TextValue "E32" '  +1.00%  '
Set-TextValue "D33" '0.105'
Set-TextValue "E33" '  +0.87%  '
Set-TextValue "D34" '57.89'
Set-TextValue "E34" '  +1.76%  '
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  -0.08%  '
Set-TextValue "D36" '3.706.31'
Set-TextValue "E36" '  -0.65%  '
Set-TextValue "D37" '3.29'
Set-TextValue "E37" '  -2.18%  '
Set-TextValue "E38" '  +3.91%  '
Set-TextValue "D39" '34.64'
Set-TextValue "E39" '  +2.56%  '
Set-TextValue "E40" '  +1.42%  '
Set-TextValue "E41" '  -0.96%  '
Set-TextValue "B42" 'TheGraph'
Set-TextValue "C42" 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue "D42" '0.336'
Set-TextValue "E42" '  +0.54%  '
Set-TextValue "B43" 'PEPE'
Set-TextValue "C43" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D43" '0.0₃0672'
Set-TextValue "E43" '  +0.66%  '
Set-TextValue "E44" '  +0.01%  '
Set-TextValue "D45" '0.0412'
Set-TextValue "E45" '  +1.04%  '
Set-TextValue "E46" '  +2.23%  '
Set-TextValue "E47" '  +0.52%  '
Set-TextValue "E48" '  -0.03%  '
Set-TextValue "E49" '  -0.45%  '
Set-TextValue "D50" '131.22'
Set-TextValue "E50" '  +3.01%  '
Set-TextValue "D51" '2.55'
Set-TextValue "E51" '  -0.90%  '
